$p = $ppt.ActivePresentation

# --- Slide 1: title placeholder ("ctrTitle") moves down slightly ---
# target a:off y goes from 1069102 EMU to 1077894 EMU (x unchanged at 790469 EMU).
# Values below are expressed in points (1 pt = 12700 EMU) and picked so that
# the host's internal float32 cast still lands on the exact target EMU.
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item(1)
$titleShape.Top = 84.87358093261719

# --- Slide 1: supervisor text box ("Google Shape;90;p13") ---
# " " -> " C"  (space run right after "Akkamahadevi")
# "Professor / Associate Professor / Assistant Professor" -> "Assistant Professor"
$supShape = $s1.Shapes.Item(4)
$supRange = $supShape.TextFrame.TextRange
$supRange.Characters(44, 1).Text = " C"
$supRange.Characters(47, 53).Text = "Assistant Professor"

# --- Slide 5: requirements heading ---
# "Software and Hardware Requirements:" -> "Software  Requirements:"
$s5 = $p.Slides.Item(5)
$reqShape = $s5.Shapes.Item(2)
$reqRange = $reqShape.TextFrame.TextRange
$reqRange.Characters(1, 35).Text = "Software  Requirements:"

# --- Slide 7: "table" picture repositioned ---
# a:off goes from x=1453877,y=1758462 to x=1441971,y=1617785 (EMU).
$s7 = $p.Slides.Item(7)
$picShape = $s7.Shapes.Item(3)
$picShape.Left = 113.54106140136719
$picShape.Top = 127.38468170166016
